$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 14,6
$data[0,0] = 46050
$data[0,1] = 12042.7213773065
$data[0,2] = 11451.5067855586
$data[0,3] = 18219.86
$data[0,4] = 7677.85681049757
$data[0,5] = 37.8959831690054
$data[1,0] = 46051
$data[1,1] = 11818.9049477622
$data[1,2] = 11471.6565526613
$data[1,3] = 12075.86
$data[1,4] = 7932.54055656622
$data[1,5] = 305.347379551146
$data[2,0] = 46052
$data[2,1] = 11987.7115914832
$data[2,2] = 11070.9881165453
$data[2,3] = 12075.86
$data[2,4] = 8071.17770724575
$data[2,5] = 294.429409324628
$data[3,0] = 46053
$data[3,1] = 4889.70006964604
$data[3,2] = 7987.07593152685
$data[3,3] = 12075.86
$data[3,4] = 8004.21437351375
$data[3,5] = 163.142929376692
$data[4,0] = 46054
$data[4,1] = 5110.2819978697
$data[4,2] = 7826.45944573543
$data[4,3] = 9791.86
$data[4,4] = 7995.24919339478
$data[4,5] = 251.243693297092
$data[5,0] = 46055
$data[5,1] = 11455.6405730248
$data[5,2] = 10928.087340475
$data[5,3] = 9791.86
$data[5,4] = 7815.5107887778
$data[5,5] = 372.989088718865
$data[6,0] = 46056
$data[6,1] = 11455.6405730248
$data[6,2] = 10595.1390077298
$data[6,3] = 9791.86
$data[6,4] = 7815.5107887778
$data[6,5] = 359.116241521151
$data[7,0] = 46057
$data[7,1] = 11455.6405730248
$data[7,2] = 10507.789924467
$data[7,3] = 9791.86
$data[7,4] = 7815.5107887778
$data[7,5] = 355.476696385201
$data[8,0] = 46058
$data[8,1] = 11455.6405730248
$data[8,2] = 10659.7928068437
$data[8,3] = 9791.86
$data[8,4] = 7815.5107887778
$data[8,5] = 361.810149817562
$data[9,0] = 46059
$data[9,1] = 11455.6405730248
$data[9,2] = 9998.88128975423
$data[9,3] = 9791.86
$data[9,4] = 7815.5107887778
$data[9,5] = 334.272169938835
$data[10,0] = 46060
$data[10,1] = 4820.25116453669
$data[10,2] = 7169.61151040991
$data[10,3] = 9791.86
$data[10,4] = 7442.15878952607
$data[10,5] = 200.829595830666
$data[11,0] = 46061
$data[11,1] = 4720.83339186724
$data[11,2] = 7213.25365095635
$data[11,3] = 9791.86
$data[11,4] = 7433.48546707663
$data[11,5] = 202.286629918041
$data[12,0] = 46062
$data[12,1] = 11066.9810198305
$data[12,2] = 10488.8662489866
$data[12,3] = 9791.86
$data[12,4] = 7339.1665303127
$data[12,5] = 334.840532470803
$data[13,0] = 46063
$data[13,1] = 11066.9810198305
$data[13,2] = 10752.605655129
$data[13,3] = 9791.86
$data[13,4] = 7339.1665303127
$data[13,5] = 345.829674393406

$ws.Range("A2:F15").Value = $data
